{"js": "const replacements = [\n  [\"2025-07-10 Thursday\", \"2025-07-11 Friday\"],\n  [\"121\u00f72=60, 1\", \"392\u00f73=130, 2\"],\n  [\"867\u00f79=96, 3\", \"855\u00f77=122, 1\"],\n  [\"326\u00f78=40, 6\", \"974\u00f75=194, 4\"],\n  [\"199\u00f72=99, 1\", \"512\u00f78=64, 0\"],\n  [\"115\u00f77=16, 3\", \"166\u00f73=55, 1\"],\n  [\"684\u00f79=76, 0\", \"122\u00f74=30, 2\"],\n  [\"890\u00f73=296, 2\", \"672\u00f74=168, 0\"],\n  [\"602\u00f77=86, 0\", \"975\u00f77=139, 2\"],\n  [\"963\u00f75=192, 3\", \"694\u00f73=231, 1\"],\n  [\"321\u00f78=40, 1\", \"701\u00f77=100, 1\"],\n  [\"228\u00f75=45, 3\", \"372\u00f78=46, 4\"],\n  [\"352\u00f78=44, 0\", \"903\u00f77=129, 0\"],\n  [\"603\u00f79=67, 0\", \"564\u00f75=112, 4\"],\n  [\"350\u00f79=38, 8\", \"988\u00f78=123, 4\"],\n  [\"827\u00f77=118, 1\", \"250\u00f75=50, 0\"],\n  [\"545\u00f76=90, 5\", \"524\u00f73=174, 2\"],\n  [\"551\u00f74=137, 3\", \"362\u00f78=45, 2\"],\n  [\"610\u00f76=101, 4\", \"309\u00f72=154, 1\"],\n  [\"140\u00f75=28, 0\", \"723\u00f76=120, 3\"],\n  [\"579\u00f75=115, 4\", \"514\u00f78=64, 2\"],\n  [\"592\u00f79=65, 7\", \"163\u00f73=54, 1\"],\n  [\"163\u00f72=81, 1\", \"905\u00f76=150, 5\"],\n  [\"911\u00f78=113, 7\", \"696\u00f75=139, 1\"],\n  [\"594\u00f78=74, 2\", \"918\u00f78=114, 6\"],\n  [\"281\u00f79=31, 2\", \"447\u00f78=55, 7\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load('items/text');\nawait context.sync();\n\nconst used = new Array(replacements.length).fill(false);\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  for (let j = 0; j < replacements.length; j++) {\n    if (used[j]) continue;\n    const [oldText, newText] = replacements[j];\n    if (text === oldText) {\n      const rng = para.getRange();\n      rng.insertText(newText, Word.InsertLocation.replace);\n      used[j] = true;\n      break;\n    }\n  }\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2025-07-10 Thursday\", \"2025-07-11 Friday\")\n    ,@(\"121\u00f72=60, 1\", \"392\u00f73=130, 2\")\n    ,@(\"867\u00f79=96, 3\", \"855\u00f77=122, 1\")\n    ,@(\"326\u00f78=40, 6\", \"974\u00f75=194, 4\")\n    ,@(\"199\u00f72=99, 1\", \"512\u00f78=64, 0\")\n    ,@(\"115\u00f77=16, 3\", \"166\u00f73=55, 1\")\n    ,@(\"684\u00f79=76, 0\", \"122\u00f74=30, 2\")\n    ,@(\"890\u00f73=296, 2\", \"672\u00f74=168, 0\")\n    ,@(\"602\u00f77=86, 0\", \"975\u00f77=139, 2\")\n    ,@(\"963\u00f75=192, 3\", \"694\u00f73=231, 1\")\n    ,@(\"321\u00f78=40, 1\", \"701\u00f77=100, 1\")\n    ,@(\"228\u00f75=45, 3\", \"372\u00f78=46, 4\")\n    ,@(\"352\u00f78=44, 0\", \"903\u00f77=129, 0\")\n    ,@(\"603\u00f79=67, 0\", \"564\u00f75=112, 4\")\n    ,@(\"350\u00f79=38, 8\", \"988\u00f78=123, 4\")\n    ,@(\"827\u00f77=118, 1\", \"250\u00f75=50, 0\")\n    ,@(\"545\u00f76=90, 5\", \"524\u00f73=174, 2\")\n    ,@(\"551\u00f74=137, 3\", \"362\u00f78=45, 2\")\n    ,@(\"610\u00f76=101, 4\", \"309\u00f72=154, 1\")\n    ,@(\"140\u00f75=28, 0\", \"723\u00f76=120, 3\")\n    ,@(\"579\u00f75=115, 4\", \"514\u00f78=64, 2\")\n    ,@(\"592\u00f79=65, 7\", \"163\u00f73=54, 1\")\n    ,@(\"163\u00f72=81, 1\", \"905\u00f76=150, 5\")\n    ,@(\"911\u00f78=113, 7\", \"696\u00f75=139, 1\")\n    ,@(\"594\u00f78=74, 2\", \"918\u00f78=114, 6\")\n    ,@(\"281\u00f79=31, 2\", \"447\u00f78=55, 7\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}"}
